$wb = $excel.ActiveWorkbook

# Rename the existing sheet "pythoned" to "Men"
$men = $wb.Worksheets.Item(1)
$men.Name = "Men"

# Add a new empty sheet "Women" right after "Men"
$women = $wb.Worksheets.Add($null, $men)
$women.Name = "Women"

# Update view state on Men: scroll position + active cell selection
$men.Activate()
$men.Range("G205").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 154
$excel.ActiveWindow.ScrollColumn = 1

# Make Women the active (selected) sheet/tab
$women.Activate()
$women.Range("A1").Select() | Out-Null
